# Bench Shield Ladder Program upload.
# - Refresh the cached "today" date field (2022/12/29 -> 2023/1/15) on the
#   slide master and every slide layout.
# - Renumber the I/O pin labels on slide 3's Digital Inputs / Digital
#   Outputs tables from %IX100.x/%QX100.x to %IX0.x/%QX0.x.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Date placeholder refresh (slide master + all custom layouts).
# ---------------------------------------------------------------------
function Update-CachedDate($shapes, $oldText, $newText) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldText) {
                $sh.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

$master = $p.SlideMaster
Update-CachedDate $master.Shapes "2022/12/29" "2023/1/15"

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-CachedDate $layouts.Item($i).Shapes "2022/12/29" "2023/1/15"
}

# ---------------------------------------------------------------------
# 2. Bench shield pin renumbering on slide 3.
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

$ixTbl = $s3.Shapes.Item(1).Table
$ixTbl.Cell(2, 2).Shape.TextFrame.TextRange.Text = "%IX0.0"
$ixTbl.Cell(3, 2).Shape.TextFrame.TextRange.Text = "%IX0.1"
$ixTbl.Cell(4, 2).Shape.TextFrame.TextRange.Text = "%IX0.2"
$ixTbl.Cell(5, 2).Shape.TextFrame.TextRange.Text = "%IX0.3"
$ixTbl.Cell(6, 2).Shape.TextFrame.TextRange.Text = "%IX0.4"

$qxTbl = $s3.Shapes.Item(2).Table
$qxTbl.Cell(2, 2).Shape.TextFrame.TextRange.Text = "%QX0.0"
$qxTbl.Cell(3, 2).Shape.TextFrame.TextRange.Text = "%QX0.1"
$qxTbl.Cell(4, 2).Shape.TextFrame.TextRange.Text = "%QX0.2"
$qxTbl.Cell(5, 2).Shape.TextFrame.TextRange.Text = "%QX0.3"

Write-Host "Applied Bench Shield Ladder Program upload edits."
